$wb = $excel.ActiveWorkbook

# Update cached numeric values across multiple sheets to reflect
# refreshed market-board price data (scheduled runner update).


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2332.8572
$ws.Range("J2").Value = 2100
$ws.Range("L2").Value = 2100
$ws.Range("N2").Value = -2326

$ws.Range("H18").Value = 8900.5
$ws.Range("I18").Value = 8900.5
$ws.Range("K18").Value = 8900.5
$ws.Range("M18").Value = -8616.5

$ws.Range("H70").Value = 21913.445
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 75000
$ws.Range("N70").Value = -75540

$ws.Range("H73").Value = 21913.445
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 75000
$ws.Range("N73").Value = -76872

$ws.Range("H129").Value = 2647.4375
$ws.Range("I129").Value = 991.8
$ws.Range("J129").Value = 3400
$ws.Range("K129").Value = 2975.4
$ws.Range("L129").Value = 10200
$ws.Range("M129").Value = 2024.6
$ws.Range("N129").Value = -20200

$ws.Range("H135").Value = 2986.0557
$ws.Range("I135").Value = 2482.3845
$ws.Range("J135").Value = 4295.6
$ws.Range("K135").Value = 22341.4605
$ws.Range("L135").Value = 38660.4
$ws.Range("M135").Value = -19806.4605
$ws.Range("N135").Value = -43730.4

$ws.Range("H137").Value = 3409.8333
$ws.Range("I137").Value = 2681.65
$ws.Range("J137").Value = 7050.75
$ws.Range("K137").Value = 8044.950000000001
$ws.Range("L137").Value = 21152.25
$ws.Range("M137").Value = -5494.950000000001
$ws.Range("N137").Value = -26252.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13293.841
$ws.Range("I32").Value = 11902.738
$ws.Range("K32").Value = 11902.738
$ws.Range("M32").Value = -11615.738

$ws.Range("H61").Value = 1756
$ws.Range("I61").Value = 1756
$ws.Range("K61").Value = 1756
$ws.Range("M61").Value = -1544

$ws.Range("H122").Value = 2847.4546
$ws.Range("I122").Value = 2883.3
$ws.Range("K122").Value = 8649.900000000001
$ws.Range("M122").Value = -6199.900000000001

$ws.Range("H126").Value = 5158.385
$ws.Range("I126").Value = 5158.385
$ws.Range("K126").Value = 15475.155
$ws.Range("M126").Value = -13005.155

$ws.Range("H132").Value = 3355.7273
$ws.Range("I132").Value = 3351.2188
$ws.Range("K132").Value = 10053.6564
$ws.Range("M132").Value = -7523.6564

$ws.Range("H136").Value = 1756
$ws.Range("I136").Value = 1756
$ws.Range("K136").Value = 5268
$ws.Range("M136").Value = -2718


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15716.353
$ws.Range("I99").Value = 16143.125
$ws.Range("K99").Value = 16143.125
$ws.Range("M99").Value = -14645.125

$ws.Range("H128").Value = 3949.3333
$ws.Range("I128").Value = 3949.3333
$ws.Range("K128").Value = 11847.9999
$ws.Range("M128").Value = -9357.999899999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 53627.45
$ws.Range("I31").Value = 65033.438
$ws.Range("J31").Value = 8003.5
$ws.Range("K31").Value = 65033.438
$ws.Range("L31").Value = 8003.5
$ws.Range("M31").Value = -64738.438
$ws.Range("N31").Value = -8593.5

$ws.Range("H34").Value = 53627.45
$ws.Range("I34").Value = 65033.438
$ws.Range("J34").Value = 8003.5
$ws.Range("K34").Value = 65033.438
$ws.Range("L34").Value = 8003.5
$ws.Range("M34").Value = -64831.438
$ws.Range("N34").Value = -8407.5

$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

$ws.Range("H58").Value = 2666
$ws.Range("I58").Value = 2999
$ws.Range("K58").Value = 2999
$ws.Range("M58").Value = -2796

$ws.Range("H69").Value = 549
$ws.Range("I69").Value = 549
$ws.Range("K69").Value = 549
$ws.Range("M69").Value = 200

$ws.Range("H72").Value = 549
$ws.Range("I72").Value = 549
$ws.Range("K72").Value = 1647
$ws.Range("M72").Value = 2097

$ws.Range("H76").Value = 5166.6665
$ws.Range("I76").Value = 5166.6665
$ws.Range("K76").Value = 5166.6665
$ws.Range("M76").Value = -4851.6665

$ws.Range("H79").Value = 5166.6665
$ws.Range("I79").Value = 5166.6665
$ws.Range("K79").Value = 5166.6665
$ws.Range("M79").Value = -4074.6665

$ws.Range("H136").Value = 2666
$ws.Range("I136").Value = 2999
$ws.Range("K136").Value = 8997
$ws.Range("M136").Value = -6447


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2258.2666
$ws.Range("I113").Value = 2154.25
$ws.Range("K113").Value = 2154.25
$ws.Range("M113").Value = 15.75

$ws.Range("H132").Value = 10312.838
$ws.Range("I132").Value = 9599.107
$ws.Range("K132").Value = 28797.321
$ws.Range("M132").Value = -26267.321


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3559
$ws.Range("I22").Value = 2338.5
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 2338.5
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -2043.5
$ws.Range("N22").Value = -6590

$ws.Range("H27").Value = 3559
$ws.Range("I27").Value = 2338.5
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 2338.5
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -2231.5
$ws.Range("N27").Value = -6214

$ws.Range("H40").Value = 18599.1
$ws.Range("I40").Value = 19777.334
$ws.Range("K40").Value = 19777.334
$ws.Range("M40").Value = -19641.334

$ws.Range("H55").Value = 2627.5898
$ws.Range("I55").Value = 1596.8695
$ws.Range("K55").Value = 1596.8695
$ws.Range("M55").Value = -1423.8695

$ws.Range("H132").Value = 5588
$ws.Range("I132").Value = 4579.4
$ws.Range("J132").Value = 6848.75
$ws.Range("K132").Value = 13738.2
$ws.Range("L132").Value = 20546.25
$ws.Range("M132").Value = -11208.2
$ws.Range("N132").Value = -25606.25

$ws.Range("H136").Value = 8010
$ws.Range("I136").Value = 7830.963
$ws.Range("K136").Value = 23492.889
$ws.Range("M136").Value = -20942.889


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1720.2
$ws.Range("I81").Value = 1720.2
$ws.Range("K81").Value = 3440.4
$ws.Range("M81").Value = -2379.4

$ws.Range("H84").Value = 1720.2
$ws.Range("I84").Value = 1720.2
$ws.Range("K84").Value = 17202
$ws.Range("M84").Value = -11898

$ws.Range("H132").Value = 6832.7617
$ws.Range("I132").Value = 6892.6
$ws.Range("K132").Value = 20677.8
$ws.Range("M132").Value = -18147.8

$ws.Range("H136").Value = 1682.2727
$ws.Range("I136").Value = 1705
$ws.Range("K136").Value = 5115
$ws.Range("M136").Value = -2565
